$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column E: values ---
$ws.Range("E3").Value = $null
$ws.Range("E4").Value = 2020
$ws.Range("E5").Value = 11.5
$ws.Range("E6").Value = 2.6
$ws.Range("E7").Value = 2
$ws.Range("E8").Value = 0.3

# --- Copy formatting from column D into column E, row by row ---
$ws.Range("D3").Copy()
$ws.Range("E3").PasteSpecial(-4122)

$ws.Range("D4").Copy()
$ws.Range("E4").PasteSpecial(-4122)

$ws.Range("D5").Copy()
$ws.Range("E5").PasteSpecial(-4122)

$ws.Range("D6").Copy()
$ws.Range("E6").PasteSpecial(-4122)

$ws.Range("D7").Copy()
$ws.Range("E7").PasteSpecial(-4122)

$ws.Range("D8").Copy()
$ws.Range("E8").PasteSpecial(-4122)

# E7 gets its own one-decimal number format (new style, distinct from D7's)
$ws.Range("E7").NumberFormat = "0.0"

# --- Selection moves to B15, matching the authored sheetView ---
$ws.Range("B15").Select()
